$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 62.75
$ws.Range("I6").Value = 61.4
$ws.Range("J6").Value = 69.5
$ws.Range("K6").Value = 184.2
$ws.Range("L6").Value = 208.5
$ws.Range("M6").Value = -72.19999999999999
$ws.Range("N6").Value = -432.5
$ws.Range("H7").Value = 5125.25
$ws.Range("I7").Value = 8000
$ws.Range("K7").Value = 8000
$ws.Range("M7").Value = -7888
$ws.Range("H14").Value = 5125.25
$ws.Range("I14").Value = 8000
$ws.Range("K14").Value = 8000
$ws.Range("M14").Value = -7809
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("H58").Value = 1266.6
$ws.Range("I58").Value = 209.71428
$ws.Range("J58").Value = 1835.6923
$ws.Range("K58").Value = 629.14284
$ws.Range("L58").Value = 5507.0769
$ws.Range("M58").Value = -479.14284
$ws.Range("N58").Value = -5807.0769
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").ClearContents()
$ws.Range("H86").Value = 598
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("H87").Value = 41795.668
$ws.Range("J87").Value = 72924.664
$ws.Range("L87").Value = 72924.664
$ws.Range("N87").Value = -75420.664
$ws.Range("H89").Value = 598
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("H90").Value = 41795.668
$ws.Range("J90").Value = 72924.664
$ws.Range("L90").Value = 218773.992
$ws.Range("N90").Value = -231253.992
$ws.Range("H116").Value = 4275.154
$ws.Range("I116").Value = 3610.889
$ws.Range("K116").Value = 3610.889
$ws.Range("M116").Value = -168.8890000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2769.0557
$ws.Range("I2").Value = 1903.8334
$ws.Range("J2").Value = 4499.5
$ws.Range("K2").Value = 1903.8334
$ws.Range("L2").Value = 4499.5
$ws.Range("M2").Value = -1790.8334
$ws.Range("N2").Value = -4725.5
$ws.Range("H44").Value = 12775.6875
$ws.Range("J44").Value = 12775.6875
$ws.Range("L44").Value = 12775.6875
$ws.Range("N44").Value = -13751.6875
$ws.Range("H45").Value = 2681.4546
$ws.Range("I45").Value = 2209.3684
$ws.Range("K45").Value = 2209.3684
$ws.Range("M45").Value = -1832.3684
$ws.Range("H55").Value = 8000
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
$ws.Range("H63").Value = 5985
$ws.Range("I63").Value = 1833
$ws.Range("J63").Value = 12213
$ws.Range("K63").Value = 1833
$ws.Range("L63").Value = 12213
$ws.Range("M63").Value = -1147
$ws.Range("N63").Value = -13585
$ws.Range("H66").Value = 5985
$ws.Range("I66").Value = 1833
$ws.Range("J66").Value = 12213
$ws.Range("K66").Value = 9165
$ws.Range("L66").Value = 61065
$ws.Range("M66").Value = -5733
$ws.Range("N66").Value = -67929
$ws.Range("H74").Value = 1853
$ws.Range("I74").Value = 1853
$ws.Range("K74").Value = 1853
$ws.Range("M74").Value = -979
$ws.Range("H77").Value = 1853
$ws.Range("I77").Value = 1853
$ws.Range("K77").Value = 9265
$ws.Range("M77").Value = -4897
$ws.Range("H116").Value = 2769.0557
$ws.Range("I116").Value = 1903.8334
$ws.Range("J116").Value = 4499.5
$ws.Range("K116").Value = 1903.8334
$ws.Range("L116").Value = 4499.5
$ws.Range("M116").Value = 390.1666
$ws.Range("N116").Value = -9087.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2769.0557
$ws.Range("I3").Value = 1903.8334
$ws.Range("J3").Value = 4499.5
$ws.Range("K3").Value = 1903.8334
$ws.Range("L3").Value = 4499.5
$ws.Range("M3").Value = -1789.8334
$ws.Range("N3").Value = -4727.5
$ws.Range("H64").Value = 511.375
$ws.Range("I64").Value = 273.25
$ws.Range("K64").Value = 273.25
$ws.Range("M64").Value = -48.25
$ws.Range("H67").Value = 511.375
$ws.Range("I67").Value = 273.25
$ws.Range("K67").Value = 273.25
$ws.Range("M67").Value = 506.75
$ws.Range("H107").Value = 3143.3438
$ws.Range("I107").Value = 1358.8
$ws.Range("J107").Value = 6117.5835
$ws.Range("K107").Value = 1358.8
$ws.Range("L107").Value = 6117.5835
$ws.Range("M107").Value = 561.2
$ws.Range("N107").Value = -9957.583500000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 92.2
$ws.Range("J7").Value = 86
$ws.Range("L7").Value = 86
$ws.Range("N7").Value = -312
$ws.Range("H31").Value = 5512.75
$ws.Range("I31").Value = 2423.261
$ws.Range("J31").Value = 9692.647000000001
$ws.Range("K31").Value = 2423.261
$ws.Range("L31").Value = 9692.647000000001
$ws.Range("M31").Value = -2128.261
$ws.Range("N31").Value = -10282.647
$ws.Range("H34").Value = 5512.75
$ws.Range("I34").Value = 2423.261
$ws.Range("J34").Value = 9692.647000000001
$ws.Range("K34").Value = 2423.261
$ws.Range("L34").Value = 9692.647000000001
$ws.Range("M34").Value = -2221.261
$ws.Range("N34").Value = -10096.647
$ws.Range("H58").Value = 2880.1177
$ws.Range("I58").Value = 2880.1177
$ws.Range("K58").Value = 2880.1177
$ws.Range("M58").Value = -2677.1177
$ws.Range("H86").Value = 7500
$ws.Range("I86").Value = 7500
$ws.Range("K86").Value = 7500
$ws.Range("M86").Value = -6377
$ws.Range("H89").Value = 7500
$ws.Range("I89").Value = 7500
$ws.Range("K89").Value = 37500
$ws.Range("M89").Value = -31884
$ws.Range("H95").Value = 20358.625
$ws.Range("J95").Value = 20358.625
$ws.Range("L95").Value = 20358.625
$ws.Range("N95").Value = -25850.625
$ws.Range("H136").Value = 2880.1177
$ws.Range("I136").Value = 2880.1177
$ws.Range("K136").Value = 8640.3531
$ws.Range("M136").Value = -6090.3531

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 104.4
$ws.Range("I7").Value = 104.666664
$ws.Range("K7").Value = 313.999992
$ws.Range("M7").Value = -201.999992
$ws.Range("H34").Value = 3299.6
$ws.Range("I34").Value = 99.5
$ws.Range("J34").Value = 5433
$ws.Range("K34").Value = 298.5
$ws.Range("L34").Value = 16299
$ws.Range("M34").Value = -214.5
$ws.Range("N34").Value = -16467
$ws.Range("H92").Value = 7149.8335
$ws.Range("J92").Value = 5974.75
$ws.Range("L92").Value = 17924.25
$ws.Range("N92").Value = -20420.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 7692307.5
$ws.Range("I11").Value = 7047619
$ws.Range("J11").Value = 10400000
$ws.Range("K11").Value = 7047619
$ws.Range("L11").Value = 10400000
$ws.Range("M11").Value = -7047480
$ws.Range("N11").Value = -10400278
$ws.Range("H14").Value = 2761.5652
$ws.Range("I14").Value = 182.1579
$ws.Range("J14").Value = 15013.75
$ws.Range("K14").Value = 182.1579
$ws.Range("L14").Value = 15013.75
$ws.Range("M14").Value = -14.15790000000001
$ws.Range("N14").Value = -15349.75
$ws.Range("H80").Value = 3984.5
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 3984.5
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7532.4443
$ws.Range("I7").Value = 6123.75
$ws.Range("J7").Value = 8659.4
$ws.Range("K7").Value = 6123.75
$ws.Range("L7").Value = 8659.4
$ws.Range("M7").Value = -6011.75
$ws.Range("N7").Value = -8883.4
$ws.Range("H122").Value = 3995.6667
$ws.Range("I122").Value = 3995.6667
$ws.Range("K122").Value = 11987.0001
$ws.Range("M122").Value = -9537.000100000001
$ws.Range("H126").Value = 7532.4443
$ws.Range("I126").Value = 6123.75
$ws.Range("J126").Value = 8659.4
$ws.Range("K126").Value = 18371.25
$ws.Range("L126").Value = 25978.2
$ws.Range("M126").Value = -15901.25
$ws.Range("N126").Value = -30918.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 497.5
$ws.Range("I6").Value = 497.5
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 497.5
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -382.5
$ws.Range("N6").ClearContents()
$ws.Range("H54").Value = 69189.60000000001
$ws.Range("J54").Value = 69189.60000000001
$ws.Range("L54").Value = 69189.60000000001
$ws.Range("N54").Value = -70229.60000000001
$ws.Range("H80").Value = 20000
$ws.Range("J80").Value = 20000
$ws.Range("L80").Value = 20000
$ws.Range("N80").Value = -21996
$ws.Range("H81").Value = 1066.6666
$ws.Range("I81").Value = 800
$ws.Range("J81").Value = 1200
$ws.Range("K81").Value = 1600
$ws.Range("L81").Value = 2400
$ws.Range("M81").Value = -539
$ws.Range("N81").Value = -4522
$ws.Range("H83").Value = 20000
$ws.Range("J83").Value = 20000
$ws.Range("L83").Value = 60000
$ws.Range("N83").Value = -69984
$ws.Range("H84").Value = 1066.6666
$ws.Range("I84").Value = 800
$ws.Range("J84").Value = 1200
$ws.Range("K84").Value = 8000
$ws.Range("L84").Value = 12000
$ws.Range("M84").Value = -2696
$ws.Range("N84").Value = -22608
$ws.Range("H135").Value = 37522.445
$ws.Range("J135").Value = 37522.445
$ws.Range("L135").Value = 37522.445
$ws.Range("N135").Value = -47662.445
